$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Value would otherwise be auto-converted to a number by Excel, so
        # enter it the way Excel's quote-prefix ('-leading) text entry does
        # to force it to stay text, then strip the resulting quote-prefix
        # formatting so the cell keeps the default (unstyled) look.
        $cell.Value = "'" + $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

# Rows 6-9: Name / Amount / Price
for ($row = 6; $row -le 9; $row++) {
    Set-TextCell $row 1 "Name"
    Set-TextCell $row 2 "Amount"
    Set-TextCell $row 3 "Price"
}

# Rows 1-4: test / Amount / Price
for ($row = 1; $row -le 4; $row++) {
    Set-TextCell $row 1 "test"
    Set-TextCell $row 2 "Amount"
    Set-TextCell $row 3 "Price"
}

# Row 5: first typed as Cola / 1L / 5eur ...
Set-TextCell 5 1 "Cola"
Set-TextCell 5 2 "1L"
Set-TextCell 5 3 "5eur"
# ... then corrected back to "test" in column A
Set-TextCell 5 1 "test"

# Row 10: Race / 5 / 0eur
Set-TextCell 10 1 "Race"
Set-TextCell 10 2 "5"
Set-TextCell 10 3 "0eur"
